$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: collapse the run of blank paragraphs + line breaks that sits
# right after the "NB: XXX si riferisce ... XXX)." paragraph down to a
# single paragraph holding one <w:br/>.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("XXX).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $markStart = $rng1.End + 1   # first of the 3 blank-paragraph marks
    $brStart = $markStart + 3    # first of the 6 <w:br/> run breaks

    # Remove the first 5 of the 6 line breaks, keep the 6th.
    $d.Range($brStart, $brStart + 5).Delete()

    # Remove the 3 blank paragraph marks that precede them.
    $d.Range($markStart, $markStart + 3).Delete()
}

# ---------------------------------------------------------------------
# Change 2: merge the two runs split by a <w:lastRenderedPageBreak/> so
# "modalità di " and "funzionamento (man/auto), farà una " become one
# run again (re-typing the same text drops the stale page-break marker).
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "modalità di funzionamento (man/auto), farà una ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "modalità di funzionamento (man/auto), farà una ", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 3: delete the extra blank "ind left=720" paragraph right before
# the document's final blank paragraph.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("sottoscrivono il dato topic.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    # $rng3.End sits on this paragraph's own ending mark; the blank
    # "ind left=720" paragraph's mark is the next character.
    $d.Range($rng3.End + 1, $rng3.End + 2).Delete()
}
